# "Generate Report for Archive"
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the 9370d91d-...md source file, on every sheet
# that reports it:
#   - Overview sheet: columns E ("zh-cn") and F ("de-de"), row 2
#   - zh-cn sheet:     column C ("Status"), row 2
#   - de-de sheet:     column C ("Status"), row 2
#
# Because the new status text is shorter than the old one, the
# generating tool also re-autosized the "Status" column(s) on every
# sheet to fit the new (narrower) content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Column width (in "characters") that yields the narrower, content-fit
# column width seen after the status text shrank.
$newStatusColWidth = 12.5

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
